$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated per-epoch validation accuracy values (column B) after freezing the encoder.
$accuracyUpdates = @{
    3 = 0.421875
    4 = 0.40625
    5 = 0.34375
    6 = 0.328125
    7 = 0.265625
    8 = 0.265625
    10 = 0.296875
    11 = 0.296875
    12 = 0.3125
    13 = 0.28125
    14 = 0.3125
    15 = 0.359375
    16 = 0.34375
    17 = 0.21875
    18 = 0.234375
    19 = 0.25
    20 = 0.25
    21 = 0.234375
    22 = 0.21875
    23 = 0.21875
    25 = 0.203125
    26 = 0.21875
    27 = 0.21875
    28 = 0.21875
    29 = 0.21875
    30 = 0.21875
    31 = 0.21875
    32 = 0.21875
    33 = 0.21875
    34 = 0.21875
    35 = 0.21875
    36 = 0.21875
    37 = 0.21875
    38 = 0.21875
    39 = 0.21875
    40 = 0.21875
    41 = 0.21875
    42 = 0.21875
    43 = 0.21875
    44 = 0.21875
    45 = 0.21875
    46 = 0.21875
    47 = 0.21875
    48 = 0.21875
    49 = 0.21875
    50 = 0.21875
    51 = 0.203125
    52 = 0.203125
    53 = 0.203125
    54 = 0.203125
    55 = 0.203125
    56 = 0.203125
    57 = 0.203125
    58 = 0.203125
    59 = 0.203125
    60 = 0.203125
    61 = 0.203125
    62 = 0.203125
    63 = 0.203125
    64 = 0.203125
    65 = 0.203125
    66 = 0.203125
    67 = 0.203125
    68 = 0.203125
    69 = 0.203125
    70 = 0.203125
    71 = 0.203125
    72 = 0.203125
    73 = 0.203125
    74 = 0.203125
    75 = 0.203125
    76 = 0.203125
    77 = 0.203125
    78 = 0.203125
    79 = 0.203125
    80 = 0.203125
    81 = 0.203125
    82 = 0.203125
    83 = 0.203125
    84 = 0.203125
    85 = 0.203125
    86 = 0.203125
    87 = 0.203125
    88 = 0.203125
    89 = 0.203125
    90 = 0.203125
    91 = 0.203125
    92 = 0.203125
    93 = 0.203125
    94 = 0.203125
    95 = 0.203125
    96 = 0.203125
    97 = 0.203125
    98 = 0.203125
    99 = 0.203125
    100 = 0.203125
    101 = 0.203125
    102 = 0.203125
    103 = 0.09375
    104 = 0.171875
    105 = 0.296875
    107 = 0.328125
    108 = 0.28125
    109 = 0.203125
    110 = 0.234375
    112 = 0.1875
    113 = 0.21875
    114 = 0.171875
    115 = 0.25
    117 = 0.234375
}

foreach ($row in $accuracyUpdates.Keys) {
    $ws.Cells.Item([int]$row, 2).Value = $accuracyUpdates[$row]
}

# Refresh the repr string stamped in column A (object id changes each run).
$newRepr = "<__main__.DisplayOutputs object at 0x7f2498306b50>"
for ($row = 102; $row -le 118; $row++) {
    $ws.Cells.Item($row, 1).Value = $newRepr
}
